$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("C2").Value = 0.1294396932959941
$ws.Range("D2").Value = 85.97371648230884
$ws.Range("C3").Value = 0.2076535069114193
$ws.Range("D3").Value = 84.85637628780276
$ws.Range("E3").Value = 0.1865958124825144
$ws.Range("C4").Value = 0.3315484433631548
$ws.Range("D4").Value = 83.08644862420654
$ws.Range("E4").Value = 0.4821737323030834
$ws.Range("C5").Value = 0.5253553379585646
$ws.Range("D5").Value = 80.31777870141499
$ws.Range("E5").Value = 0.9445416094092752
$ws.Range("C6").Value = 0.8224605018333396
$ws.Range("D6").Value = 76.07341921748962
$ws.Range("E6").Value = 1.65334964322481
$ws.Range("C7").Value = 1.263411897844627
$ws.Range("D7").Value = 69.7741135601855
$ws.Range("E7").Value = 2.705333687994596
$ws.Range("C8").Value = 1.885020863059941
$ws.Range("D8").Value = 60.89398548568101
$ws.Range("E8").Value = 4.188315076436845
$ws.Range("C9").Value = 2.693529737658339
$ws.Range("D9").Value = 49.34385870570391
$ws.Range("E9").Value = 6.117186248693023
$ws.Range("C10").Value = 3.624791330147982
$ws.Range("D10").Value = 36.04012167013759
$ws.Range("E10").Value = 8.338910333632597
$ws.Range("C11").Value = 4.529107093435749
$ws.Range("D11").Value = 23.12132505174092
$ws.Range("E11").Value = 10.49634936890484
$ws.Range("C12").Value = 5.240658897767408
$ws.Range("D12").Value = 12.95629927557438
$ws.Range("E12").Value = 12.19390867352466
$ws.Range("C13").Value = 5.693150829998912
$ws.Range("D13").Value = 6.492128815124324
$ws.Range("E13").Value = 13.27342514041982
$ws.Range("C14").Value = 5.935920177477171
$ws.Range("D14").Value = 3.023995279720622
$ws.Range("E14").Value = 13.85260344083223
$ws.Range("C15").Value = 6.052834252301664
$ws.Range("D15").Value = 1.353794210799295
$ws.Range("E15").Value = 14.1315270193421
$ws.Range("C16").Value = 6.105981495155647
$ws.Range("D16").Value = 0.5945478843138013
$ws.Range("E16").Value = 14.25832115586518
$ws.Range("C17").Value = 6.129481412492744
$ws.Range("D17").Value = 0.2588347794981508
$ws.Range("E17").Value = 14.31438524436939
$ws.Range("C18").Value = 6.139742511593304
$ws.Range("D18").Value = 0.1122476494901516
$ws.Range("E18").Value = 14.33886529508072
$ws.Range("C19").Value = 6.144198149853797
$ws.Range("D19").Value = 0.04859567434024974
$ws.Range("E19").Value = 14.34949517493076
$ws.Range("C20").Value = 6.146128222095857
$ws.Range("D20").Value = 0.02102321373937761
$ws.Range("E20").Value = 14.3540997758511
$ws.Range("C21").Value = 6.146963402637169
$ws.Range("D21").Value = 0.009092063149217719
$ws.Range("E21").Value = 14.35609227799966
